# TP_2or_fold3.xlsx — "Add files via upload"
#
# The source workbook is a 1-sheet transition-probability matrix. The
# uploaded revision renames the "ScreenRecStarted" state to "0_unstated"
# everywhere it's used (the column header in G1, and the four row labels
# in A27:A30 that are built from that state name), and leaves the cursor
# selection on E14 instead of G9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the (former) "ScreenRecStarted" column.
$ws.Range("G1").Value = "0_unstated"

# Row labels referencing the renamed state (transition pairs).
$ws.Range("A27").Value = "0_unstated1_Scanning"
$ws.Range("A28").Value = "0_unstated3_Reading"
$ws.Range("A29").Value = "0_unstated5_Unknown "
$ws.Range("A30").Value = "0_unstated0_unstated"

# Match the selection saved in the uploaded file.
[void]$ws.Range("E14").Select()
